# Apply updated dSF (column F) values after repulling/recalculating data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 0
    6  = -4
    7  = -2
    9  = 3
    17 = -3
    21 = 1
    22 = -2
    28 = 7
    31 = -3
    34 = -3
    43 = -4
    48 = 1
    56 = -5
    57 = -2
    63 = -5
    64 = -5
    66 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
